# Wiring guide edit: fix PDP port numbers, flag the Victor SPX device-type
# label, and rename a couple of variable names.
#
# Summary of the real edits made to the workbook (per the authoritative diff):
#   PDP sheet, row 12 (cargo pickup wheels / Victor SPX):
#       - Device Type label flagged with a red "(needs fix)" suffix
#       - PDP port changed from 6 to 11
#   PDP sheet, row 13 (hatch pickup wheels / Victor SPX):
#       - Device Type label flagged with a red "(needs fix)" suffix
#       - PDP port changed from 11 to 6
#       - Variable Name renamed from "intake" to "intakeMotor"
#   PCM sheet, row 2 (lower/raise hatch pickup):
#       - Variable Name renamed from "pickupSol" to "rotatorSol"
#   View-state bookkeeping (active cell selections) updated to match
#   where the author was last working.

$wb = $excel.ActiveWorkbook

$wsPDP = $wb.Worksheets.Item("PDP")
$wsPCM = $wb.Worksheets.Item("PCM")

# --- Swap the PDP breaker-port numbers between the two Victor SPX rows ---
$wsPDP.Range("E12").Value = 11
$wsPDP.Range("E13").Value = 6

# --- Rename variables ---
$wsPDP.Range("J13").Value = "intakeMotor"
$wsPCM.Range("F2").Value = "rotatorSol"

# --- Mark the "Victor SPX" device type as needing a fix, with the ---
# --- "(needs fix)" portion in bold red text.                      ---
$wsPDP.Range("B12").Value = "Victor SPX (needs fix)"
$charsB12 = $wsPDP.Range("B12").Characters(12, 11)
$charsB12.Font.Bold = $true
$charsB12.Font.Color = 255

# Copy the same rich-text label down to B13 so both cells share the
# identical formatted string (rather than creating a duplicate entry).
$wsPDP.Range("B12").Copy()
$wsPDP.Range("B13").PasteSpecial()
$excel.CutCopyMode = $false

# --- Restore the author's final selection state: last looked at PCM!F3, ---
# --- then ended on PDP!B13 (PDP remains the active/visible tab).        ---
$wsPCM.Select()
$wsPCM.Range("F3").Select()

$wsPDP.Range("B13").Select()
